# Fix 2 sites, add 5 cars
# (This workbook's slice of the change: fix alyans-auto.ru / zd-auto.ru links
# and prices, and drop one stale alyans-auto.ru offer.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - Changan CS35 Plus: alyans-auto.ru link refreshed
$ws.Range("K10").Value = "https://alyans-auto.ru/auto/auto_17299.html"

# Row 20 - Changan UNI-K: alyans-auto.ru link refreshed
$ws.Range("K20").Value = "https://alyans-auto.ru/auto/auto_16980.html"

# Row 21 - Changan UNI-S: alyans-auto.ru link refreshed (both min_price_url and site column)
$ws.Range("E21").Value = "https://alyans-auto.ru/auto/auto_18043.html"
$ws.Range("K21").Value = "https://alyans-auto.ru/auto/auto_18043.html"

# Row 23 - Changan UNI-V: stale alyans-auto.ru offer removed entirely
$ws.Range("J23:K23").ClearContents()

# Row 44 - EXEED LX: zd-auto.ru price update
$ws.Range("L44").Value = 2739900

# Row 59 - Geely Atlas Pro: alyans-auto.ru link refreshed
$ws.Range("K59").Value = "https://alyans-auto.ru/auto/auto_17094.html"

# Row 61 - Geely Cityray: price + alyans-auto.ru link refreshed
$ws.Range("D61").Value = 3214190
$ws.Range("E61").Value = "https://alyans-auto.ru/auto/auto_19149.html"
$ws.Range("J61").Value = 3214190
$ws.Range("K61").Value = "https://alyans-auto.ru/auto/auto_19149.html"

# Row 68 - Geely Monjaro: alyans-auto.ru link refreshed
$ws.Range("K68").Value = "https://alyans-auto.ru/auto/auto_18909.html"

# Row 73 - Great Wall POER: alyans-auto.ru price update
$ws.Range("J73").Value = 3499000

# Row 74 - Haval Dargo: alyans-auto.ru price + link refreshed
$ws.Range("J74").Value = 3649000
$ws.Range("K74").Value = "https://alyans-auto.ru/auto/auto_18190.html"

# Row 76 - Haval F7: alyans-auto.ru price + link refreshed
$ws.Range("J76").Value = 2849000
$ws.Range("K76").Value = "https://alyans-auto.ru/auto/auto_19315.html"

# Row 85 - Haval Jolion: alyans-auto.ru price + link refreshed
$ws.Range("J85").Value = 2149000
$ws.Range("K85").Value = "https://alyans-auto.ru/auto/auto_19275.html"

# Row 87 - Haval M6: alyans-auto.ru price update
$ws.Range("J87").Value = 2279000

# Row 102 - JAECOO J7: price update (min_price + alyans-auto.ru price)
$ws.Range("D102").Value = 2389900
$ws.Range("J102").Value = 2389900

# Row 103 - JAECOO J8: alyans-auto.ru price + link refreshed
$ws.Range("J103").Value = 4259000
$ws.Range("K103").Value = "https://alyans-auto.ru/auto/auto_17152.html"

# Row 128 - Knewstar 001: price + alyans-auto.ru link refreshed
$ws.Range("D128").Value = 4393190
$ws.Range("E128").Value = "https://alyans-auto.ru/auto/auto_19305.html"
$ws.Range("J128").Value = 4393190
$ws.Range("K128").Value = "https://alyans-auto.ru/auto/auto_19305.html"

# Row 150 - Lada Vesta Sedan: alyans-auto.ru price + link refreshed
$ws.Range("J150").Value = 1575000
$ws.Range("K150").Value = "https://alyans-auto.ru/auto/auto_17708.html"

# Row 173 - OMODA S5: zd-auto.ru link fixed
$ws.Range("M173").Value = "https://zd-auto.ru/catalog/omoda/omoda_s5"
